$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New multi-line note for row 50 (B50) describing work done
$note = "Improved sample code for saving file in any directory" + [char]10 + "Created Use Case diagram" + [char]10 + "Drafted Activity diagram"

$ws.Range("B50").Value = $note
$ws.Range("D50").Value = 5
$ws.Range("D51").Value = 4
$ws.Rows.Item(50).RowHeight = 57.75

# Update the selection to reflect where the user left off editing
$ws.Range("C49").Select()
